# Add "Molar Vol (L/mol)" and "M. Vol. err (L/mol)" columns to the GSA
# (Gas Sorption Analysis) results table, right after the "Fug. err (MPa)"
# column and before "Mass frac." — i.e. insert two new columns at N:O,
# pushing the existing Mass frac. / Mass frac. err / Dual Mode Pred (CC/CC)
# / Dual Mode Pred Err (CC/CC) columns two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column N (shifts N:Q -> P:S).
$ws.Range("N1:O1").EntireColumn.Insert()

# Fill in the headers for the two newly inserted columns on the header row.
$ws.Range("N11").Value = "Molar Vol (L/mol)"
$ws.Range("O11").Value = "M. Vol. err (L/mol)"
